$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: management_address changed (switch was re-IP'd)
$ws.Range('F2').Value = '10.0.0.20'

# Row 3: new ICX6450 switch 1
$ws.Range('A3').Value = '00:06:cc:4e:24:c1_0/10/1/3'
$ws.Range('B3').Value = 'ICX6450_template_A'
$ws.Range('C3').Value = 'ZTP_ICX6450_SW1'
$ws.Range('D3').Value = 'The Lab'
$ws.Range('E3').Value = 'Tim Braly'
$ws.Range('F3').Value = '10.0.0.3'
$ws.Range('G3').Value = '255.255.255.0'
$ws.Range('H3').Value = '10.0.0.1'
$ws.Range('I3').Value = 10
$ws.Range('J3').Value = 'e 1/1/1'
$ws.Range('K3').Value = 'e 1/1/2'
$ws.Range('L3').Value = 201
$ws.Range('M3').Value = 'e 1/1/3 to 1/1/44'
$ws.Range('N3').Value = 1001
$ws.Range('O3').Value = 'e 1/1/3 to 1/1/6'
$ws.Range('P3').Value = 2001
$ws.Range('Q3').Value = 'e 1/1/7 to 1/1/12'

# Row 4: new ICX6450 switch 2
$ws.Range('A4').Value = '00:06:cc:4e:24:c1_0/10/1/5'
$ws.Range('B4').Value = 'ICX6450_template_A'
$ws.Range('C4').Value = 'ZTP_ICX6450_SW2'
$ws.Range('D4').Value = 'The Lab'
$ws.Range('E4').Value = 'Tim Braly'
$ws.Range('F4').Value = '10.0.0.4'
$ws.Range('G4').Value = '255.255.255.0'
$ws.Range('H4').Value = '10.0.0.1'
$ws.Range('I4').Value = 10
$ws.Range('J4').Value = 'e 1/1/1'
$ws.Range('K4').Value = 'e 1/1/2'
$ws.Range('L4').Value = 201
$ws.Range('M4').Value = 'e 1/1/3 to 1/1/44'
$ws.Range('N4').Value = 1001
$ws.Range('O4').Value = 'e 1/1/3 to 1/1/6'
$ws.Range('P4').Value = 2001
$ws.Range('Q4').Value = 'e 1/1/7 to 1/1/12'

# Row 5: device serial number lookup row
$ws.Range('A5').Value = 'CTG2549K09L'
$ws.Range('B5').Value = 'ICX6450_template_A'
$ws.Range('C5').Value = 'ZTP_ICX6450_SERIAL'
$ws.Range('D5').Value = 'The Lab'
$ws.Range('E5').Value = 'Tim Braly'
$ws.Range('F5').Value = '10.0.0.2'
$ws.Range('G5').Value = '255.255.255.0'
$ws.Range('H5').Value = '10.0.0.1'
$ws.Range('I5').Value = 10
$ws.Range('J5').Value = 'e 1/1/1'
$ws.Range('K5').Value = 'e 1/1/2'
$ws.Range('L5').Value = 201
$ws.Range('M5').Value = 'e 1/1/3 to 1/1/44'
$ws.Range('N5').Value = 1001
$ws.Range('O5').Value = 'e 1/1/3 to 1/1/6'
$ws.Range('P5').Value = 2001
$ws.Range('Q5').Value = 'e 1/1/7 to 1/1/12'

# Move selection to A5 (also clears the stale topLeftCell scroll position)
$ws.Range('A5').Select()
